# Apply the row permutation described by the diff: the 14 data rows
# (rows 2-15) of the "Artfynd" sheet get reordered into a new sequence.
# Column C (Valideringsstatus) and columns I, K, P, S, T, U, V, W, Y, Z, AA, AB,
# AD, AE, AG, AT, AW, AX, AY are identical across all 14 rows, so they are
# left untouched; only A, B, D, E, F, G, H, M, Q, R need to be rewritten per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 = old row 7 (Id 96780135)
$ws.Range("A2").Value = 96780135
$ws.Range("B2").Value = 89392
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = "Ullticka"
$ws.Range("G2").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H2").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q2").Value = 658733.0048414731
$ws.Range("R2").Value = 6636920.163120084

# New row 3 = old row 8 (Id 96779183)
$ws.Range("A3").Value = 96779183
$ws.Range("B3").Value = 98520
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("Q3").Value = 658704.2164550385
$ws.Range("R3").Value = 6637062.857129692

# New row 4 = old row 9 (Id 96780122)
$ws.Range("A4").Value = 96780122
$ws.Range("B4").Value = 89832
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 1209
$ws.Range("F4").Value = "Rynkskinn"
$ws.Range("G4").Value = "Phlebia centrifuga"
$ws.Range("H4").Value = "P.Karst."
$ws.Range("Q4").Value = 658733.0048414731
$ws.Range("R4").Value = 6636920.163120084

# New row 5 = old row 10 (Id 96779825)
$ws.Range("A5").Value = 96779825
$ws.Range("B5").Value = 89376
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4660
$ws.Range("F5").Value = "Rävticka"
$ws.Range("G5").Value = "Inocutis rheades"
$ws.Range("H5").Value = "(Pers.) Fiasson & Niemelä"
$ws.Range("Q5").Value = 658637.0597997338
$ws.Range("R5").Value = 6636982.990721731

# New row 6 = old row 2 (Id 96775736)
$ws.Range("A6").Value = 96775736
$ws.Range("B6").Value = 89412
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 5442
$ws.Range("F6").Value = "Tallticka"
$ws.Range("G6").Value = "Porodaedalea pini"
$ws.Range("H6").Value = "(Brot.) Murrill"
$ws.Range("Q6").Value = 658730.5226168972
$ws.Range("R6").Value = 6637449.43415721

# New row 7 = old row 3 (Id 96777744)
$ws.Range("A7").Value = 96777744
$ws.Range("B7").Value = 98520
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 222498
$ws.Range("F7").Value = "Blåsippa"
$ws.Range("G7").Value = "Hepatica nobilis"
$ws.Range("H7").Value = "Schreb."
$ws.Range("Q7").Value = 658769.3765012868
$ws.Range("R7").Value = 6637283.535847809

# New row 8 = old row 11 (Id 96779798)
$ws.Range("A8").Value = 96779798
$ws.Range("B8").Value = 43464
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 101735
$ws.Range("F8").Value = "Jättesvampmal"
$ws.Range("G8").Value = "Scardia boletella"
$ws.Range("H8").Value = "(Fabricius, 1794)"
$ws.Range("Q8").Value = 658637.0597997338
$ws.Range("R8").Value = 6636982.990721731
$ws.Range("M8").Value = "äldre gnagspår"

# New row 9 = old row 12 (Id 96780278)
$ws.Range("A9").Value = 96780278
$ws.Range("B9").Value = 93132
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 2671
$ws.Range("F9").Value = "Fällmossa"
$ws.Range("G9").Value = "Antitrichia curtipendula"
$ws.Range("H9").Value = "(Hedw.) Brid."
$ws.Range("Q9").Value = 658722.9088558007
$ws.Range("R9").Value = 6636991.191442309

# New row 10 = old row 13 (Id 96779993)
$ws.Range("A10").Value = 96779993
$ws.Range("B10").Value = 43464
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 101735
$ws.Range("F10").Value = "Jättesvampmal"
$ws.Range("G10").Value = "Scardia boletella"
$ws.Range("H10").Value = "(Fabricius, 1794)"
$ws.Range("Q10").Value = 658739.4120713262
$ws.Range("R10").Value = 6636888.229354058
$ws.Range("M10").Value = "äldre gnagspår"

# New row 11 = old row 4 (Id 96778360)
$ws.Range("A11").Value = 96778360
$ws.Range("B11").Value = 98520
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 222498
$ws.Range("F11").Value = "Blåsippa"
$ws.Range("G11").Value = "Hepatica nobilis"
$ws.Range("H11").Value = "Schreb."
$ws.Range("Q11").Value = 658743.2312543363
$ws.Range("R11").Value = 6637305.564015599
$ws.Range("M11").Value = ""

# New row 12 = old row 14 (Id 96780357)
$ws.Range("A12").Value = 96780357
$ws.Range("B12").Value = 98520
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 222498
$ws.Range("F12").Value = "Blåsippa"
$ws.Range("G12").Value = "Hepatica nobilis"
$ws.Range("H12").Value = "Schreb."
$ws.Range("Q12").Value = 658747.5451754113
$ws.Range("R12").Value = 6637110.504147635

# New row 13 = old row 15 (Id 96780175)
$ws.Range("A13").Value = 96780175
$ws.Range("B13").Value = 90005
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 1339
$ws.Range("F13").Value = "Brandticka"
$ws.Range("G13").Value = "Pycnoporellus fulgens"
$ws.Range("H13").Value = "(Fr.) Donk"
$ws.Range("Q13").Value = 658742.8409314866
$ws.Range("R13").Value = 6636937.694258579
$ws.Range("M13").Value = ""

# New row 14 = old row 5 (Id 96777066)
$ws.Range("A14").Value = 96777066
$ws.Range("B14").Value = 89392
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = "Ullticka"
$ws.Range("G14").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H14").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q14").Value = 658712.8441804391
$ws.Range("R14").Value = 6637414.964914286

# New row 15 = old row 6 (Id 96780597)
$ws.Range("A15").Value = 96780597
$ws.Range("B15").Value = 89392
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 1202
$ws.Range("F15").Value = "Ullticka"
$ws.Range("G15").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H15").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q15").Value = 658693.8994370478
$ws.Range("R15").Value = 6637469.000542388
